$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSA2LGXFJX%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T220646Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIFtWC8cnPZK3rrnW15Xud4KGj90FJHU4lbfn8vm5lNX2AiEA3OYAwEOTcRDfWXOFS1s7YliaYlo7GuwaZLuYbMpYr7sq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJn4pdO5znUq%2BifOsiq6A%2BWScnmDdZeoDcaZdwIsUg8qIo%2FtN9XK84Ldhv%2Fciptm4qXYfNn6rRQzf74pZUn9c2emYcnb9o35QTlr3e7im%2FJA2SdYoOme4ua846hQ35YOHV2d9Yj%2BpAFsJYNOiXclF2ic6zHIPQ3ohMOplwgmEZKQbW8vqOPssU5vyuf%2FMLB9OuIw%2Fmu3qiHbI2th%2FaHxyakrie6BxhQPDKo6sFlyhEW0SiuXIP0h2qx%2Bu6XVXkk8mX5JMxh89HIddC2%2FB%2BXk7Ei%2Fu9F%2FEHcXd76u4EBxiTqAx%2F9prxjSLarpvDbmaX9Mlgrj960NuUyKKgsqBOs3dAoPTv2bIVXIGA%2Fy6QiHVcICRXFmv9n6vWkY%2FKEmZI3thLvfCNnD31Eck9K5P6aUy6hCUsi%2FHYKkWPxNrnb8O2S82Vq7R%2B3wPA%2BsJDdLxOvZsv0he40sUeT53Luersyly%2FARGKTeMS8%2FZAeon5ecSom43j%2FE2OLP7BA1Amx3dG77GV%2BNwXKi%2FLPbW8%2Ffllqf4t4OMlD8TZJsf39qxvHjMA7fB2yWVfn3ba5a0RPvPxOfRCNPtz5nMS0a9uu%2BU5lyXFN%2F8tehUHO4%2BMQw8t22vwY6pQG1qyYlpdgLzqW552L0kh%2FFRD3mIca5QXsolm%2FBma27jakabLIepOg1zGB%2BQ%2F1z2xmVILUy0ykFi45QnGjmcjF7xkq0TQIWC4oaswqrT2%2FShNct6uWKU3feKJeDbopyLpSkqEG9BsdDJ9Eii3RJUtv1TOzBBwLsCmupJml4MS0UERpygP7qM%2BI%2FCE4S7UL06Xwuo0zgAwSkveiM%2ByI0AfZ%2Ba1yr3Cg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2006%3A06%3A46%20GMT&X-Amz-Signature=1270e4a0f1793f3ae4b1f5435ef5fe5d428806df10abdeebaacd2f4f27c9315d"

$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSA2LGXFJX%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T220646Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIFtWC8cnPZK3rrnW15Xud4KGj90FJHU4lbfn8vm5lNX2AiEA3OYAwEOTcRDfWXOFS1s7YliaYlo7GuwaZLuYbMpYr7sq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJn4pdO5znUq%2BifOsiq6A%2BWScnmDdZeoDcaZdwIsUg8qIo%2FtN9XK84Ldhv%2Fciptm4qXYfNn6rRQzf74pZUn9c2emYcnb9o35QTlr3e7im%2FJA2SdYoOme4ua846hQ35YOHV2d9Yj%2BpAFsJYNOiXclF2ic6zHIPQ3ohMOplwgmEZKQbW8vqOPssU5vyuf%2FMLB9OuIw%2Fmu3qiHbI2th%2FaHxyakrie6BxhQPDKo6sFlyhEW0SiuXIP0h2qx%2Bu6XVXkk8mX5JMxh89HIddC2%2FB%2BXk7Ei%2Fu9F%2FEHcXd76u4EBxiTqAx%2F9prxjSLarpvDbmaX9Mlgrj960NuUyKKgsqBOs3dAoPTv2bIVXIGA%2Fy6QiHVcICRXFmv9n6vWkY%2FKEmZI3thLvfCNnD31Eck9K5P6aUy6hCUsi%2FHYKkWPxNrnb8O2S82Vq7R%2B3wPA%2BsJDdLxOvZsv0he40sUeT53Luersyly%2FARGKTeMS8%2FZAeon5ecSom43j%2FE2OLP7BA1Amx3dG77GV%2BNwXKi%2FLPbW8%2Ffllqf4t4OMlD8TZJsf39qxvHjMA7fB2yWVfn3ba5a0RPvPxOfRCNPtz5nMS0a9uu%2BU5lyXFN%2F8tehUHO4%2BMQw8t22vwY6pQG1qyYlpdgLzqW552L0kh%2FFRD3mIca5QXsolm%2FBma27jakabLIepOg1zGB%2BQ%2F1z2xmVILUy0ykFi45QnGjmcjF7xkq0TQIWC4oaswqrT2%2FShNct6uWKU3feKJeDbopyLpSkqEG9BsdDJ9Eii3RJUtv1TOzBBwLsCmupJml4MS0UERpygP7qM%2BI%2FCE4S7UL06Xwuo0zgAwSkveiM%2ByI0AfZ%2Ba1yr3Cg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2006%3A06%3A46%20GMT&X-Amz-Signature=9b7798f2e465eb530901d28e2e99de1f58d6ca9054999da7acc98345c6e99979"

$ws.Range("K4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSA2LGXFJX%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T220646Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIFtWC8cnPZK3rrnW15Xud4KGj90FJHU4lbfn8vm5lNX2AiEA3OYAwEOTcRDfWXOFS1s7YliaYlo7GuwaZLuYbMpYr7sq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJn4pdO5znUq%2BifOsiq6A%2BWScnmDdZeoDcaZdwIsUg8qIo%2FtN9XK84Ldhv%2Fciptm4qXYfNn6rRQzf74pZUn9c2emYcnb9o35QTlr3e7im%2FJA2SdYoOme4ua846hQ35YOHV2d9Yj%2BpAFsJYNOiXclF2ic6zHIPQ3ohMOplwgmEZKQbW8vqOPssU5vyuf%2FMLB9OuIw%2Fmu3qiHbI2th%2FaHxyakrie6BxhQPDKo6sFlyhEW0SiuXIP0h2qx%2Bu6XVXkk8mX5JMxh89HIddC2%2FB%2BXk7Ei%2Fu9F%2FEHcXd76u4EBxiTqAx%2F9prxjSLarpvDbmaX9Mlgrj960NuUyKKgsqBOs3dAoPTv2bIVXIGA%2Fy6QiHVcICRXFmv9n6vWkY%2FKEmZI3thLvfCNnD31Eck9K5P6aUy6hCUsi%2FHYKkWPxNrnb8O2S82Vq7R%2B3wPA%2BsJDdLxOvZsv0he40sUeT53Luersyly%2FARGKTeMS8%2FZAeon5ecSom43j%2FE2OLP7BA1Amx3dG77GV%2BNwXKi%2FLPbW8%2Ffllqf4t4OMlD8TZJsf39qxvHjMA7fB2yWVfn3ba5a0RPvPxOfRCNPtz5nMS0a9uu%2BU5lyXFN%2F8tehUHO4%2BMQw8t22vwY6pQG1qyYlpdgLzqW552L0kh%2FFRD3mIca5QXsolm%2FBma27jakabLIepOg1zGB%2BQ%2F1z2xmVILUy0ykFi45QnGjmcjF7xkq0TQIWC4oaswqrT2%2FShNct6uWKU3feKJeDbopyLpSkqEG9BsdDJ9Eii3RJUtv1TOzBBwLsCmupJml4MS0UERpygP7qM%2BI%2FCE4S7UL06Xwuo0zgAwSkveiM%2ByI0AfZ%2Ba1yr3Cg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2006%3A06%3A46%20GMT&X-Amz-Signature=825a17c9988eeeb40458ee576d52c82e41b3a3cc5c82bf04edca47b8c8abeeb9"

$ws.Range("K5").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSA2LGXFJX%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T220646Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIFtWC8cnPZK3rrnW15Xud4KGj90FJHU4lbfn8vm5lNX2AiEA3OYAwEOTcRDfWXOFS1s7YliaYlo7GuwaZLuYbMpYr7sq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJn4pdO5znUq%2BifOsiq6A%2BWScnmDdZeoDcaZdwIsUg8qIo%2FtN9XK84Ldhv%2Fciptm4qXYfNn6rRQzf74pZUn9c2emYcnb9o35QTlr3e7im%2FJA2SdYoOme4ua846hQ35YOHV2d9Yj%2BpAFsJYNOiXclF2ic6zHIPQ3ohMOplwgmEZKQbW8vqOPssU5vyuf%2FMLB9OuIw%2Fmu3qiHbI2th%2FaHxyakrie6BxhQPDKo6sFlyhEW0SiuXIP0h2qx%2Bu6XVXkk8mX5JMxh89HIddC2%2FB%2BXk7Ei%2Fu9F%2FEHcXd76u4EBxiTqAx%2F9prxjSLarpvDbmaX9Mlgrj960NuUyKKgsqBOs3dAoPTv2bIVXIGA%2Fy6QiHVcICRXFmv9n6vWkY%2FKEmZI3thLvfCNnD31Eck9K5P6aUy6hCUsi%2FHYKkWPxNrnb8O2S82Vq7R%2B3wPA%2BsJDdLxOvZsv0he40sUeT53Luersyly%2FARGKTeMS8%2FZAeon5ecSom43j%2FE2OLP7BA1Amx3dG77GV%2BNwXKi%2FLPbW8%2Ffllqf4t4OMlD8TZJsf39qxvHjMA7fB2yWVfn3ba5a0RPvPxOfRCNPtz5nMS0a9uu%2BU5lyXFN%2F8tehUHO4%2BMQw8t22vwY6pQG1qyYlpdgLzqW552L0kh%2FFRD3mIca5QXsolm%2FBma27jakabLIepOg1zGB%2BQ%2F1z2xmVILUy0ykFi45QnGjmcjF7xkq0TQIWC4oaswqrT2%2FShNct6uWKU3feKJeDbopyLpSkqEG9BsdDJ9Eii3RJUtv1TOzBBwLsCmupJml4MS0UERpygP7qM%2BI%2FCE4S7UL06Xwuo0zgAwSkveiM%2ByI0AfZ%2Ba1yr3Cg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2006%3A06%3A46%20GMT&X-Amz-Signature=4fb5838c2b71c0d7c528c54ba60d43bb50778170caecb3e5a571612d971a9ed0"

$ws.Range("K6").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSA2LGXFJX%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T220646Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIFtWC8cnPZK3rrnW15Xud4KGj90FJHU4lbfn8vm5lNX2AiEA3OYAwEOTcRDfWXOFS1s7YliaYlo7GuwaZLuYbMpYr7sq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJn4pdO5znUq%2BifOsiq6A%2BWScnmDdZeoDcaZdwIsUg8qIo%2FtN9XK84Ldhv%2Fciptm4qXYfNn6rRQzf74pZUn9c2emYcnb9o35QTlr3e7im%2FJA2SdYoOme4ua846hQ35YOHV2d9Yj%2BpAFsJYNOiXclF2ic6zHIPQ3ohMOplwgmEZKQbW8vqOPssU5vyuf%2FMLB9OuIw%2Fmu3qiHbI2th%2FaHxyakrie6BxhQPDKo6sFlyhEW0SiuXIP0h2qx%2Bu6XVXkk8mX5JMxh89HIddC2%2FB%2BXk7Ei%2Fu9F%2FEHcXd76u4EBxiTqAx%2F9prxjSLarpvDbmaX9Mlgrj960NuUyKKgsqBOs3dAoPTv2bIVXIGA%2Fy6QiHVcICRXFmv9n6vWkY%2FKEmZI3thLvfCNnD31Eck9K5P6aUy6hCUsi%2FHYKkWPxNrnb8O2S82Vq7R%2B3wPA%2BsJDdLxOvZsv0he40sUeT53Luersyly%2FARGKTeMS8%2FZAeon5ecSom43j%2FE2OLP7BA1Amx3dG77GV%2BNwXKi%2FLPbW8%2Ffllqf4t4OMlD8TZJsf39qxvHjMA7fB2yWVfn3ba5a0RPvPxOfRCNPtz5nMS0a9uu%2BU5lyXFN%2F8tehUHO4%2BMQw8t22vwY6pQG1qyYlpdgLzqW552L0kh%2FFRD3mIca5QXsolm%2FBma27jakabLIepOg1zGB%2BQ%2F1z2xmVILUy0ykFi45QnGjmcjF7xkq0TQIWC4oaswqrT2%2FShNct6uWKU3feKJeDbopyLpSkqEG9BsdDJ9Eii3RJUtv1TOzBBwLsCmupJml4MS0UERpygP7qM%2BI%2FCE4S7UL06Xwuo0zgAwSkveiM%2ByI0AfZ%2Ba1yr3Cg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2006%3A06%3A46%20GMT&X-Amz-Signature=76fc15e58d3a41ea658de0a66f6527fc0873e02b186ab12c9561df969d29e767"
